$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$updates = @{
    2  = -6
    3  = 0
    5  = -5
    6  = -5
    7  = -2
    8  = -3
    9  = -6
    10 = -4
    12 = -2
    13 = -6
    14 = 1
    15 = -2
    16 = 6
    17 = 1
    18 = -2
    19 = -5
    22 = -1
    23 = -6
    24 = -2
    25 = 2
    26 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
